$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f8c4340>),
                (''model'',
                 LGBMClassifier(class_weight=''balanced'', colsample_bytree=0.9,
                                learning_rate=0.05, max_depth=3,
                                min_child_samples=7, num_leaves=2,
                                random_state=42, subsample=0.5))])'
$ws.Range("B2").Value = 0.6262393162393163
$ws.Range("C2").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f633460>, ''scaler'': RobustScaler(), ''model__subsample'': 0.5, ''model__num_leaves'': 2, ''model__min_child_samples'': 7, ''model__max_depth'': 3, ''model__learning_rate'': 0.05, ''model__colsample_bytree'': 0.9, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''gbdt''}'
$ws.Range("D2").Value = 0.8034621418857302
$ws.Range("E2").Value = 0.4804183316683316
$ws.Range("F2").Value = 0.823529411764706
$ws.Range("G2").Value = 0.8883674012824031
$ws.Range("H2").Value = 0.5490119047619049
$ws.Range("I2").Value = 0.7777777777777778
$ws.Range("J2").Value = 0.7465957446808511
$ws.Range("K2").Value = 0.4516666666666667
$ws.Range("L2").Value = 0.875
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 1 1 1 1 1 0 0 1 1 1 1 0 0 1 0 0 1 1 1 1 1 1 1]'

# Row 3
$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f633520>),
                (''model'',
                 LGBMClassifier(class_weight=''balanced'', colsample_bytree=0.7,
                                learning_rate=0.05, max_depth=3,
                                min_child_samples=1, num_leaves=10,
                                random_state=42, subsample=0.7))])'
$ws.Range("B3").Value = 0.5651748251748251
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f62ebe0>, ''scaler'': RobustScaler(), ''model__subsample'': 0.7, ''model__num_leaves'': 10, ''model__min_child_samples'': 1, ''model__max_depth'': 3, ''model__learning_rate'': 0.05, ''model__colsample_bytree'': 0.7, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''gbdt''}'
$ws.Range("D3").Value = 0.7452972856481779
$ws.Range("E3").Value = 0.4840977078477079
$ws.Range("F3").Value = 0.6875
$ws.Range("G3").Value = 0.8524756293190217
$ws.Range("H3").Value = 0.6905238095238095
$ws.Range("I3").Value = 0.6875
$ws.Range("J3").Value = 0.6786170212765958
$ws.Range("K3").Value = 0.4008333333333334
$ws.Range("L3").Value = 0.6875
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[0 1 0 1 1 1 0 1 1 0 1 0 1 0 1 1 1 1 1 0 1 0 1 1]'

# Row 4
$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f633b50>),
                (''model'',
                 LGBMClassifier(class_weight=''balanced'', colsample_bytree=0.9,
                                learning_rate=0.01, max_depth=3,
                                min_child_samples=1, num_leaves=10,
                                random_state=42, subsample=0.9))])'
$ws.Range("B4").Value = 0.5881313131313132
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f62e9d0>, ''scaler'': StandardScaler(), ''model__subsample'': 0.9, ''model__num_leaves'': 10, ''model__min_child_samples'': 1, ''model__max_depth'': 3, ''model__learning_rate'': 0.01, ''model__colsample_bytree'': 0.9, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''gbdt''}'
$ws.Range("D4").Value = 0.7671851930890795
$ws.Range("E4").Value = 0.448747446997447
$ws.Range("F4").Value = 0.6451612903225806
$ws.Range("G4").Value = 0.8200559410906767
$ws.Range("H4").Value = 0.5357678571428571
$ws.Range("I4").Value = 0.8333333333333334
$ws.Range("J4").Value = 0.7431111111111111
$ws.Range("K4").Value = 0.413
$ws.Range("L4").Value = 0.5263157894736842
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 1 1 0 0 1 0 1 1 0 0 0 0 1 0 1 0 1 1 0 1 1 0 1]'

# Row 5
$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f62e040>),
                (''model'',
                 LGBMClassifier(boosting_type=''dart'', class_weight=''balanced'',
                                colsample_bytree=0.9, learning_rate=0.05,
                                max_depth=1, min_child_samples=1, num_leaves=5,
                                random_state=42, subsample=0.7))])'
$ws.Range("B5").Value = 0.583030303030303
$ws.Range("C5").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f7660d0>, ''scaler'': RobustScaler(), ''model__subsample'': 0.7, ''model__num_leaves'': 5, ''model__min_child_samples'': 1, ''model__max_depth'': 1, ''model__learning_rate'': 0.05, ''model__colsample_bytree'': 0.9, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''dart''}'
$ws.Range("D5").Value = 0.74468610914908
$ws.Range("E5").Value = 0.4582979797979798
$ws.Range("F5").Value = 0.64
$ws.Range("G5").Value = 0.8338224460027911
$ws.Range("H5").Value = 0.5855238095238094
$ws.Range("I5").Value = 0.7272727272727273
$ws.Range("J5").Value = 0.6835714285714285
$ws.Range("K5").Value = 0.4050000000000001
$ws.Range("L5").Value = 0.5714285714285714
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[0 0 1 1 1 0 0 0 1 0 1 1 1 0 0 0 1 0 0 1 0 0 1 1]'

# Row 6
$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9d316880>),
                (''model'',
                 LGBMClassifier(class_weight=''balanced'', colsample_bytree=0.5,
                                learning_rate=0.05, max_depth=1,
                                min_child_samples=1, num_leaves=10,
                                random_state=42, subsample=0.9))])'
$ws.Range("B6").Value = 0.6661383061383062
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa9f766550>, ''scaler'': MinMaxScaler(), ''model__subsample'': 0.9, ''model__num_leaves'': 10, ''model__min_child_samples'': 1, ''model__max_depth'': 1, ''model__learning_rate'': 0.05, ''model__colsample_bytree'': 0.5, ''model__class_weight'': ''balanced'', ''model__boosting_type'': ''gbdt''}'
$ws.Range("D6").Value = 0.7762729792191194
$ws.Range("E6").Value = 0.587638306138306
$ws.Range("F6").Value = 0.5384615384615385
$ws.Range("G6").Value = 0.8581702307670799
$ws.Range("H6").Value = 0.6363194444444444
$ws.Range("I6").Value = 0.4666666666666667
$ws.Range("J6").Value = 0.7123076923076923
$ws.Range("K6").Value = 0.5600000000000002
$ws.Range("L6").Value = 0.6363636363636364
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 1 1 1 0 0 0 0 1 0 1 1 0 0 0 1 0 1 1 1 1 1]'
